# Updates cryptos list values to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.154.97"
$ws.Range("E2").Value = "  +0.69%  "

$ws.Range("D3").Value = "3.937.86"
$ws.Range("E3").Value = "  +3.81%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "471.19"
$ws.Range("E5").Value = "  +9.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.01"
$ws.Range("E6").Value = "  +3.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.733"
$ws.Range("E9").Value = "  -0.41%  "

$ws.Range("E10").Value = "  +7.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000335"
$ws.Range("E11").Value = "  +6.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.40"
$ws.Range("E12").Value = "  +1.07%  "

$ws.Range("D13").Value = "4.573.80"
$ws.Range("E13").Value = "  +4.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.38"
$ws.Range("E14").Value = "  -0.85%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.26"
$ws.Range("E15").Value = "  +2.33%  "

$ws.Range("D16").Value = "3.962.87"
$ws.Range("E16").Value = "  +3.80%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.138"
$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.85"
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("E19").Value = "  +1.81%  "

$ws.Range("D20").Value = "67.520.79"
$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "437.56"
$ws.Range("E21").Value = "  +6.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.40"
$ws.Range("E22").Value = "  +4.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.50"
$ws.Range("E23").Value = "  -1.81%  "

$ws.Range("E24").Value = "  +2.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.61"
$ws.Range("E25").Value = "  +7.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "39.03"
$ws.Range("E26").Value = "  +5.99%  "

$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.76"
$ws.Range("E27").Value = "  +2.47%  "

$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.24"
$ws.Range("E28").Value = "  +4.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.77"
$ws.Range("E29").Value = "  +2.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "723.09"
$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.58"
$ws.Range("E31").Value = "  -1.50%  "

$ws.Range("E32").Value = "  -2.18%  "

$ws.Range("E33").Value = "  +4.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.83"
$ws.Range("E34").Value = "  +2.58%  "

$ws.Range("E35").Value = "  +1.37%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.85"
$ws.Range("E36").Value = "  +3.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("E38").Value = "  +15.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.36"
$ws.Range("E39").Value = "  -6.06%  "

$ws.Range("E40").Value = "  +0.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.06"
$ws.Range("E41").Value = "  +4.78%  "

$ws.Range("E42").Value = "  -6.00%  "

$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("E44").Value = "  +5.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.82"
$ws.Range("E46").Value = "  +4.53%  "

$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.48"
$ws.Range("E47").Value = "  +4.11%  "

$ws.Range("E48").Value = "  +4.50%  "

$ws.Range("E49").Value = "  +3.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.16"
$ws.Range("E50").Value = "  -2.87%  "

$ws.Range("E51").Value = "  +1.20%  "
